$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @(8.93309423163338, 8.148153864499756, 12.5708367331005, 29.12342893346305, 30.10333526103014, 14.7329636641907, 23.15210226019871, 9.514435860272796, 9.754485261953725, 59.77304688006132, 22.58507632388758)
    3 = @(9.051328042029457, 8.227226711303633, 12.61445115287429, 29.51530660825306, 30.51694268806513, 14.86362221334477, 23.43330033730727, 9.544634348360713, 9.77079496311832, 56.35621862236447, 22.83321332009884)
    4 = @(9.126803405203672, 8.277750074019325, 12.64623632790727, 29.76818166981928, 30.78911310573458, 14.94791086018614, 23.61456701396926, 9.566001921399087, 9.783710905364526, 54.13913652609128, 22.99417312831751)
    5 = @(9.158287105408176, 8.298836734785755, 12.66043011866964, 29.87429484375665, 30.90452448674829, 14.98327869392542, 23.69059033866465, 9.575412904466923, 9.789697469629031, 53.20595263876976, 23.06191318371205)
    6 = @(9.163558934847908, 8.302368294593837, 12.66286143249739, 29.89209942820638, 30.9239579539175, 14.98921300655662, 23.70334388986851, 9.577017880885776, 9.790735026397586, 53.04921302170159, 23.07329063805477)
    7 = @(9.127225058497682, 8.278032437105621, 12.64642274841545, 29.76960036445092, 30.79065146803321, 14.9483837181033, 23.61558357467132, 9.566126001765184, 9.783788722480354, 54.12667116461887, 22.99507801871776)
    8 = @(8.973265372355829, 8.175010060740711, 12.58482534959716, 29.25599221145439, 30.24211742566301, 14.77716986465094, 23.2472665851873, 9.52425759163067, 9.759502701996887, 58.61977918210538, 22.66883888841021)
    9 = @(8.694063900761225, 7.988537631555161, 12.50455442233004, 28.34696940436285, 29.31474182216143, 14.47377014838336, 22.5937540634107, 9.464896458834783, 9.735194731728592, 66.47800570293771, 22.09807234276285)
    10 = @(8.502591253743331, 7.86088150428598, 12.471398990931, 27.74025206338901, 28.72929155757689, 14.2707331766855, 22.15617919154093, 9.435605859499177, 9.731957780433959, 71.66474567775207, 21.72184172484331)
    11 = @(8.418410724996233, 7.80480704646756, 12.46216503014899, 27.47783700747057, 28.485059117939, 14.18271743546371, 21.96651584558492, 9.425492684474611, 9.733748223920397, 73.89658839914628, 21.56029593830913)
    12 = @(8.386950945078601, 7.783858049571767, 12.45952940656311, 27.38045115580831, 28.3958653640111, 14.15001687555614, 21.8960612347366, 9.422133456422188, 9.734902559786031, 74.72338469894365, 21.50052668312501)
    13 = @(8.393707839691022, 7.788357138555019, 12.46005841876436, 27.4013361221887, 28.41492653879698, 14.1570314798796, 21.91117383367515, 9.422835865259964, 9.734632654423915, 74.5461353678147, 21.51333619748298)
    14 = @(8.415814155618918, 7.803077855671272, 12.46193082355665, 27.46978504048097, 28.47765444269281, 14.18001450868774, 21.96069201186824, 9.425206840973905, 9.733833604634887, 73.96497718526585, 21.55535039401028)
    15 = @(8.429409217612365, 7.812131800850697, 12.46319045921808, 27.51197147717689, 28.51650929354571, 14.19417432829696, 21.99120182790085, 9.426720652744939, 9.733406407627788, 73.60661118860278, 21.58126895454309)
    16 = @(8.508151216058653, 7.86458611648352, 12.4721218538475, 27.75767749298921, 28.74570696258457, 14.27657301849619, 22.1687644810589, 9.436332132408854, 9.731907038821522, 71.51632055839823, 21.73259435660255)
    17 = @(8.557203280822776, 7.89727519678576, 12.47911319220721, 27.91191101066365, 28.89205053021409, 14.32823803706641, 22.28010999348311, 9.443057035454936, 9.731828014996706, 70.20129630965175, 21.82790362874411)
    18 = @(8.585691847469402, 7.91626517151348, 12.48368371987857, 28.00189855861769, 28.9782956309637, 14.35836327209901, 22.34503592062994, 9.447226754109527, 9.732089260935416, 69.4329140502511, 21.88362541591584)
    19 = @(8.595384916387863, 7.9227272139875, 12.48532495358395, 28.03258511442542, 29.00784922499054, 14.36863326700803, 22.36716989477826, 9.448690111494038, 9.732230172660095, 69.17068988266817, 21.90264614873431)
    20 = @(8.551953149903198, 7.893775937167836, 12.4783119573443, 27.89536025460466, 28.87625673351738, 14.32269586161623, 22.26816557991492, 9.442309867914638, 9.731804633496472, 70.34252649159779, 21.81766423083812)
    21 = @(8.40930968194335, 7.798746300065725, 12.46135732336385, 27.44962579426481, 28.45913945342511, 14.17324672904186, 21.94611009182468, 9.424497591382803, 9.734055320567316, 74.13617543549145, 21.54297148870851)
    22 = @(8.318516686934839, 7.738300302743749, 12.45530483933856, 27.16989957904241, 28.20577370952011, 14.07924199537758, 21.7436034419398, 9.415602287219414, 9.738306014411275, 76.50857939595542, 21.37164322512832)
    23 = @(8.366752827820937, 7.770410102065719, 12.45806843204417, 27.31812330909714, 28.33919937919181, 14.12907697465247, 21.85094957080423, 9.420095717242063, 9.735780639841735, 75.25216389195651, 21.46232588934065)
    24 = @(8.554325837214286, 7.895357339966688, 12.47867248054698, 27.90283875273307, 28.88339054473035, 14.32520016398225, 22.27356281011127, 9.442646717732085, 9.731814249579731, 70.27871481623421, 21.82229057633097)
    25 = @(8.767183469788721, 8.037332345340952, 12.52181707899993, 28.58224822644983, 29.54919732036213, 14.55236987490627, 22.76312271797652, 9.478473892018654, 9.739232424882939, 64.45493694816959, 22.24498947289688)
}

$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "L", "M", "O")

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}
